$wb = $excel.ActiveWorkbook

# --- Registration_Details sheet: update row 2 test data ---
$ws3 = $wb.Worksheets.Item("Registration_Details")
$ws3.Range("A2").Value = "pr23"
$ws3.Range("B2").Value = "pr"
$ws3.Range("C2").Value = "wb"
$ws3.Range("D2").Value = "pr23@gmail.com"
$ws3.Range("E2").Value = "Password@1234"
$ws3.Range("F2").Value = "Password@1234"

# Add hyperlinks on the password / confirm-password cells
$ws3.Hyperlinks.Add($ws3.Range("E2"), "Password@1234")
$ws3.Hyperlinks.Add($ws3.Range("F2"), "Password@1234")

# Make Registration_Details the active/selected sheet
$ws3.Select()
$ws3.Range("E4").Select()
